$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    # Force the cell to be written as text (matching the source inlineStr cells),
    # even when the value looks numeric (e.g. "0.9998"), then clear the temporary
    # number format so no stray style is left behind on the cell.
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Refreshed Price (D) / Volume(1h) (E) quotes
Set-TextValue "D2" "27.116.84"
Set-TextValue "E2" "  -2.86%  "
Set-TextValue "D3" "1.867.99"
Set-TextValue "E3" "  -2.18%  "
Set-TextValue "D4" "0.9998"
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "307.10"
Set-TextValue "E5" "  -1.95%  "
Set-TextValue "D6" "0.9996"
Set-TextValue "E6" "  +0.12%  "
Set-TextValue "D7" "0.5106"
Set-TextValue "E7" "  +1.78%  "
Set-TextValue "D8" "0.3738"
Set-TextValue "E8" "  -2.10%  "
Set-TextValue "D9" "0.07136"
Set-TextValue "E9" "  -2.44%  "
Set-TextValue "D10" "0.8881"
Set-TextValue "E10" "  -2.58%  "
Set-TextValue "D11" "20.63"
Set-TextValue "E11" "  -2.76%  "
Set-TextValue "D12" "0.07533"
Set-TextValue "E12" "  -1.77%  "
Set-TextValue "D13" "1.834.08"
Set-TextValue "E13" "  -4.54%  "
Set-TextValue "D14" "5.317"
Set-TextValue "E14" "  -3.08%  "
Set-TextValue "D15" "89.29"
Set-TextValue "E15" "  -3.99%  "
Set-TextValue "D16" "0.9977"
Set-TextValue "E16" "  -0.07%  "
Set-TextValue "D17" "0.000008456"
Set-TextValue "E17" "  -3.28%  "
Set-TextValue "D18" "14.13"
Set-TextValue "E18" "  -3.89%  "
Set-TextValue "D19" "1.000"
Set-TextValue "E19" "  +0.15%  "
Set-TextValue "D20" "27.156.25"
Set-TextValue "E20" "  -2.88%  "
Set-TextValue "D21" "5.057"
Set-TextValue "E21" "  -2.46%  "
Set-TextValue "D22" "2.100.82"
Set-TextValue "E22" "  -2.90%  "
Set-TextValue "D23" "10.57"
Set-TextValue "E23" "  -2.79%  "
Set-TextValue "D24" "6.487"
Set-TextValue "E24" "  -1.87%  "
Set-TextValue "D25" "149.67"
Set-TextValue "E25" "  -2.27%  "
Set-TextValue "E26" "  +0.15%  "
Set-TextValue "D27" "17.93"
Set-TextValue "E27" "  -2.77%  "
Set-TextValue "D28" "2.102"
Set-TextValue "E28" "  -4.92%  "
Set-TextValue "D29" "112.66"
Set-TextValue "E29" "  -2.44%  "
Set-TextValue "D30" "4.751"
Set-TextValue "E30" "  -3.70%  "
Set-TextValue "D31" "4.682"
Set-TextValue "E31" "  -3.65%  "
Set-TextValue "D32" "0.09037"
Set-TextValue "E32" "  +0.03%  "
Set-TextValue "D33" "0.05134"
Set-TextValue "E33" "  -2.93%  "
Set-TextValue "D34" "3.094"
Set-TextValue "E34" "  -3.51%  "
Set-TextValue "D35" "1.162"
Set-TextValue "E35" "  -6.06%  "
Set-TextValue "D36" "0.7377"
Set-TextValue "E36" "  -5.49%  "
Set-TextValue "D37" "0.02044"
Set-TextValue "E37" "  -2.00%  "
Set-TextValue "D38" "2.500"
Set-TextValue "E38" "  -4.41%  "
Set-TextValue "D39" "3.044"
Set-TextValue "E39" "  -0.79%  "
Set-TextValue "D40" "1.076"
Set-TextValue "E40" "  -1.54%  "
Set-TextValue "D41" "0.5320"
Set-TextValue "E41" "  -4.10%  "
Set-TextValue "D42" "6.606"
Set-TextValue "E42" "  -3.98%  "
Set-TextValue "D43" "116.46"
Set-TextValue "E43" "  +2.53%  "
Set-TextValue "D44" "8.346"
Set-TextValue "E44" "  -2.14%  "
Set-TextValue "E45" "  -3.04%  "

# Rows 46 and 47 swap content: Decentraland and PaxDollar exchange positions
Set-TextValue "B46" "Decentraland"
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.4635"
Set-TextValue "E46" "  -4.21%  "

Set-TextValue "B47" "PaxDollar"
Set-TextValue "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D47" "0.9993"
Set-TextValue "E47" "  +0.12%  "

Set-TextValue "D48" "10.00"
Set-TextValue "E48" "  -5.54%  "
Set-TextValue "D49" "1.569"
Set-TextValue "E49" "  -4.51%  "
Set-TextValue "D50" "64.49"
Set-TextValue "E50" "  -4.65%  "
Set-TextValue "D51" "36.51"
Set-TextValue "E51" "  -1.75%  "
